$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Ligand/Receptor-expressing cell counts (E, K: 1 -> 3) and all dependent
# expression/specificity metrics recomputed accordingly, per Dr Hou's advice.
$updates = @{
    2 = @{ "E"=3; "G"=24.62964533333333; "H"=73.888936; "I"=0.7460691037955391; "J"=0.746069103795539; "K"=3; "M"=64.221457; "N"=192.664371; "O"=0.8366610301096816; "P"=0.8366610301096816; "Q"=1581.751708699917; "R"=14235.76537829926; "S"=0.6242069449145827; "T"=0.6242069449145826 }
    3 = @{ "E"=3; "G"=24.62964533333333; "H"=73.888936; "I"=0.7460691037955391; "J"=0.746069103795539; "K"=3; "M"=10.582537; "N"=31.747611; "O"=0.1378666371208897; "P"=0.1378666371208896; "Q"=260.6441330368773; "R"=2345.797197331896; "S"=0.102858038400087; "T"=0.1028580384000869 }
    4 = @{ "E"=3; "G"=24.62964533333333; "H"=73.888936; "I"=0.7460691037955391; "J"=0.746069103795539; "K"=3; "M"=1.955236666666667; "N"=5.86571; "O"=0.02547233276942866; "P"=0.02547233276942866; "Q"=48.15678564272889; "R"=433.41107078456; "S"=0.01900412048086939; "T"=0.01900412048086938 }
    5 = @{ "E"=3; "G"=1.769696333333333; "H"=5.309089; "I"=0.05360677100832464; "J"=0.05360677100832464; "K"=3; "M"=64.221457; "N"=192.664371; "O"=0.8366610301096816; "P"=0.8366610301096816; "Q"=113.6524769742243; "R"=1022.872292768019; "S"=0.04485069625267871; "T"=0.04485069625267871 }
    6 = @{ "E"=3; "G"=1.769696333333333; "H"=5.309089; "I"=0.05360677100832464; "J"=0.05360677100832464; "K"=3; "M"=10.582537; "N"=31.747611; "O"=0.1378666371208897; "P"=0.1378666371208896; "Q"=18.72787692626433; "R"=168.550892336379; "S"=0.007390585245827323; "T"=0.00739058524582732 }
    7 = @{ "E"=3; "G"=1.769696333333333; "H"=5.309089; "I"=0.05360677100832464; "J"=0.05360677100832464; "K"=3; "M"=1.955236666666667; "N"=5.86571; "O"=0.02547233276942866; "P"=0.02547233276942866; "Q"=3.460175159798889; "R"=31.14157643819; "S"=0.001365489509818606; "T"=0.001365489509818606 }
    8 = @{ "E"=3; "G"=1.275024; "H"=3.825072; "I"=0.03862239996247084; "J"=0.03862239996247083; "K"=3; "M"=64.221457; "N"=192.664371; "O"=0.8366610301096816; "P"=0.8366610301096816; "Q"=81.883898989968; "R"=736.9550909097121; "S"=0.03231385693790898; "T"=0.03231385693790897 }
    9 = @{ "E"=3; "G"=1.275024; "H"=3.825072; "I"=0.03862239996247084; "J"=0.03862239996247083; "K"=3; "M"=10.582537; "N"=31.747611; "O"=0.1378666371208897; "P"=0.1378666371208896; "Q"=13.492988655888; "R"=121.436897902992; "S"=0.00532474040036383; "T"=0.005324740400363829 }
    10 = @{ "E"=3; "G"=1.275024; "H"=3.825072; "I"=0.03862239996247084; "J"=0.03862239996247083; "K"=3; "M"=1.955236666666667; "N"=5.86571; "O"=0.02547233276942866; "P"=0.02547233276942866; "Q"=2.49297367568; "R"=22.43676308112; "S"=0.0009838026241980261; "T"=0.0009838026241980259 }
    11 = @{ "E"=3; "G"=5.338186666666666; "H"=16.01456; "I"=0.1617017252336654; "J"=0.1617017252336654; "K"=3; "M"=64.221457; "N"=192.664371; "O"=0.8366610301096816; "P"=0.8366610301096816; "Q"=342.8261254713066; "R"=3085.43512924176; "S"=0.1352895320045112; "T"=0.1352895320045112 }
    12 = @{ "E"=3; "G"=5.338186666666666; "H"=16.01456; "I"=0.1617017252336654; "J"=0.1617017252336654; "K"=3; "M"=10.582537; "N"=31.747611; "O"=0.1378666371208897; "P"=0.1378666371208896; "Q"=56.49155791290666; "R"=508.42402121616; "S"=0.02229327307461156; "T"=0.02229327307461155 }
    13 = @{ "E"=3; "G"=5.338186666666666; "H"=16.01456; "I"=0.1617017252336654; "J"=0.1617017252336654; "K"=3; "M"=1.955236666666667; "N"=5.86571; "O"=0.02547233276942866; "P"=0.02547233276942866; "Q"=10.43741830417778; "R"=93.9367647376; "S"=0.004118920154542645; "T"=0.004118920154542643 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
